$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 940.3333
$ws.Range("I32").Value = 940.3333
$ws.Range("K32").Value = 940.3333
$ws.Range("M32").Value = -614.3333

$ws.Range("H51").Value = 3803.6072
$ws.Range("I51").Value = 4628.5713
$ws.Range("J51").Value = 2978.6428
$ws.Range("K51").Value = 4628.5713
$ws.Range("L51").Value = 2978.6428
$ws.Range("M51").Value = -4144.5713
$ws.Range("N51").Value = -3946.6428

$ws.Range("H74").Value = 7672.1816
$ws.Range("I74").Value = 4878.8
$ws.Range("J74").Value = 10000
$ws.Range("K74").Value = 4878.8
$ws.Range("L74").Value = 10000
$ws.Range("M74").Value = -3942.8
$ws.Range("N74").Value = -11872

$ws.Range("H77").Value = 7672.1816
$ws.Range("I77").Value = 4878.8
$ws.Range("J77").Value = 10000
$ws.Range("K77").Value = 24394
$ws.Range("L77").Value = 50000
$ws.Range("M77").Value = -19714
$ws.Range("N77").Value = -59360

$ws.Range("H86").Value = 187502140
$ws.Range("I86").Value = 100002580
$ws.Range("J86").Value = 625000000
$ws.Range("K86").Value = 100002580
$ws.Range("L86").Value = 625000000
$ws.Range("M86").Value = -100001457
$ws.Range("N86").Value = -625002246

$ws.Range("H89").Value = 187502140
$ws.Range("I89").Value = 100002580
$ws.Range("J89").Value = 625000000
$ws.Range("K89").Value = 500012900
$ws.Range("L89").Value = 3125000000
$ws.Range("M89").Value = -500007284
$ws.Range("N89").Value = -3125011232

$ws.Range("H98").Value = 1810.421
$ws.Range("I98").Value = 1602.3889
$ws.Range("K98").Value = 1602.3889
$ws.Range("M98").Value = -104.3888999999999

$ws.Range("H100").Value = 3073.4
$ws.Range("I100").Value = 1964.5
$ws.Range("K100").Value = 1964.5
$ws.Range("M100").Value = -1423.5

$ws.Range("H103").Value = 979.7646999999999
$ws.Range("J103").Value = 1096.4
$ws.Range("L103").Value = 3289.2
$ws.Range("N103").Value = -4461.200000000001

$ws.Range("H105").Value = 45887.25
$ws.Range("J105").Value = 35299.715
$ws.Range("L105").Value = 35299.715
$ws.Range("N105").Value = -42287.715

$ws.Range("H107").Value = 1762.4193
$ws.Range("I107").Value = 2073.76
$ws.Range("J107").Value = 465.16666
$ws.Range("K107").Value = 2073.76
$ws.Range("L107").Value = 465.16666
$ws.Range("M107").Value = -153.7600000000002
$ws.Range("N107").Value = -4305.16666

$ws.Range("H112").Value = 4737452
$ws.Range("J112").Value = 5263046.5
$ws.Range("L112").Value = 15789139.5
$ws.Range("N112").Value = -15791355.5

$ws.Range("H113").Value = 13173.777
$ws.Range("I113").Value = 16150.571
$ws.Range("J113").Value = 2755
$ws.Range("K113").Value = 16150.571
$ws.Range("L113").Value = 2755
$ws.Range("M113").Value = -12896.571
$ws.Range("N113").Value = -9263

$ws.Range("H122").Value = 1810.421
$ws.Range("I122").Value = 1602.3889
$ws.Range("K122").Value = 4807.1667
$ws.Range("M122").Value = -2357.1667

$ws.Range("H127").Value = 1068.2858
$ws.Range("I127").Value = 1068.2858
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 3204.8574
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = 1755.1426
$ws.Range("N127").ClearContents()

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws.Range("H132").Value = 3788.75
$ws.Range("I132").Value = 4197.024
$ws.Range("J132").Value = 930.8333
$ws.Range("K132").Value = 12591.072
$ws.Range("L132").Value = 2792.4999
$ws.Range("M132").Value = -10061.072
$ws.Range("N132").Value = -7852.4999

$ws.Range("H135").Value = 213.16667
$ws.Range("I135").Value = 189.48276
$ws.Range("K135").Value = 1705.34484
$ws.Range("M135").Value = 829.6551599999998

$ws.Range("H137").Value = 18513.568
$ws.Range("I137").Value = 28336.156
$ws.Range("J137").Value = 1970.2632
$ws.Range("K137").Value = 85008.46799999999
$ws.Range("L137").Value = 5910.7896
$ws.Range("M137").Value = -82458.46799999999
$ws.Range("N137").Value = -11010.7896

$ws.Range("H138").Value = 2311.14
$ws.Range("I138").Value = 1024.3541
$ws.Range("J138").Value = 3498.9424
$ws.Range("K138").Value = 3073.0623
$ws.Range("L138").Value = 10496.8272
$ws.Range("M138").Value = 2066.9377
$ws.Range("N138").Value = -20776.8272

$ws.Range("H140").Value = 599999.5
$ws.Range("J140").Value = 599999.5
$ws.Range("L140").Value = 599999.5
$ws.Range("N140").Value = -610359.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19003546
$ws.Range("I32").Value = 17922956
$ws.Range("J32").Value = 28574470
$ws.Range("K32").Value = 17922956
$ws.Range("L32").Value = 28574470
$ws.Range("M32").Value = -17922669
$ws.Range("N32").Value = -28575044

$ws.Range("H45").Value = 4500
$ws.Range("I45").Value = 4785.7144
$ws.Range("J45").Value = 3500
$ws.Range("K45").Value = 4785.7144
$ws.Range("L45").Value = 3500
$ws.Range("M45").Value = -4408.7144
$ws.Range("N45").Value = -4254

$ws.Range("H60").Value = 30000
$ws.Range("I60").Value = 30000
$ws.Range("K60").Value = 30000
$ws.Range("M60").Value = -29267

$ws.Range("H61").Value = 2850.0557
$ws.Range("I61").Value = 2495.8333
$ws.Range("K61").Value = 2495.8333
$ws.Range("M61").Value = -2283.8333

$ws.Range("H63").Value = 3012.1875
$ws.Range("I63").Value = 2018.9
$ws.Range("J63").Value = 4667.6665
$ws.Range("K63").Value = 2018.9
$ws.Range("L63").Value = 4667.6665
$ws.Range("M63").Value = -1332.9
$ws.Range("N63").Value = -6039.6665

$ws.Range("H66").Value = 3012.1875
$ws.Range("I66").Value = 2018.9
$ws.Range("J66").Value = 4667.6665
$ws.Range("K66").Value = 10094.5
$ws.Range("L66").Value = 23338.3325
$ws.Range("M66").Value = -6662.5
$ws.Range("N66").Value = -30202.3325

$ws.Range("H74").Value = 1922.1296
$ws.Range("I74").Value = 1887.0217
$ws.Range("K74").Value = 1887.0217
$ws.Range("M74").Value = -1013.0217

$ws.Range("H77").Value = 1922.1296
$ws.Range("I77").Value = 1887.0217
$ws.Range("K77").Value = 9435.1085
$ws.Range("M77").Value = -5067.1085

$ws.Range("H92").Value = 44011.168
$ws.Range("J92").Value = 44011.168
$ws.Range("L92").Value = 44011.168
$ws.Range("N92").Value = -49003.168

$ws.Range("H97").Value = 703.9706
$ws.Range("I97").Value = 571.7406999999999
$ws.Range("K97").Value = 571.7406999999999
$ws.Range("M97").Value = -75.74069999999995

$ws.Range("H102").Value = 1975.9375
$ws.Range("I102").Value = 1627.3077
$ws.Range("J102").Value = 3486.6667
$ws.Range("K102").Value = 1627.3077
$ws.Range("L102").Value = 3486.6667
$ws.Range("M102").Value = -5.307700000000068
$ws.Range("N102").Value = -6730.6667

$ws.Range("H110").Value = 2540
$ws.Range("I110").Value = 2443.8572
$ws.Range("J110").Value = 2764.3333
$ws.Range("K110").Value = 2443.8572
$ws.Range("L110").Value = 2764.3333
$ws.Range("M110").Value = -398.8571999999999
$ws.Range("N110").Value = -6854.3333

$ws.Range("H122").Value = 3565.25
$ws.Range("I122").Value = 2037
$ws.Range("K122").Value = 6111
$ws.Range("M122").Value = -3661

$ws.Range("H124").Value = 51566.75
$ws.Range("J124").Value = 51566.75
$ws.Range("L124").Value = 51566.75
$ws.Range("N124").Value = -61386.75

$ws.Range("H125").Value = 89338.55499999999
$ws.Range("J125").Value = 100578.14
$ws.Range("L125").Value = 100578.14
$ws.Range("N125").Value = -110418.14

$ws.Range("H127").Value = 107095.62
$ws.Range("J127").Value = 107095.62
$ws.Range("L127").Value = 107095.62
$ws.Range("N127").Value = -117015.62

$ws.Range("H133").Value = 100000
$ws.Range("J133").Value = 100000
$ws.Range("L133").Value = 100000
$ws.Range("N133").Value = -105060

$ws.Range("H136").Value = 2850.0557
$ws.Range("I136").Value = 2495.8333
$ws.Range("K136").Value = 7487.499899999999
$ws.Range("M136").Value = -4937.499899999999

$ws.Range("H141").Value = 257999.6
$ws.Range("J141").Value = 257999.6
$ws.Range("L141").Value = 257999.6
$ws.Range("N141").Value = -268359.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 916.3333
$ws.Range("I22").Value = 749.5
$ws.Range("K22").Value = 749.5
$ws.Range("M22").Value = -576.5

$ws.Range("H94").Value = 1006.5862
$ws.Range("I94").Value = 1008.5
$ws.Range("J94").Value = 1004.2308
$ws.Range("K94").Value = 1008.5
$ws.Range("L94").Value = 1004.2308
$ws.Range("M94").Value = -557.5
$ws.Range("N94").Value = -1906.2308

$ws.Range("H99").Value = 3023.15
$ws.Range("I99").Value = 2397.8
$ws.Range("J99").Value = 3648.5
$ws.Range("K99").Value = 2397.8
$ws.Range("L99").Value = 3648.5
$ws.Range("M99").Value = -899.8000000000002
$ws.Range("N99").Value = -6644.5

$ws.Range("H105").Value = 1967.5186
$ws.Range("I105").Value = 1610.8096
$ws.Range("K105").Value = 1610.8096
$ws.Range("M105").Value = 136.1904

$ws.Range("H107").Value = 1517.6666
$ws.Range("I107").Value = 917.9167
$ws.Range("J107").Value = 3916.6667
$ws.Range("K107").Value = 917.9167
$ws.Range("L107").Value = 3916.6667
$ws.Range("M107").Value = 1002.0833
$ws.Range("N107").Value = -7756.6667

$ws.Range("H125").Value = 79659.336
$ws.Range("J125").Value = 79659.336
$ws.Range("L125").Value = 79659.336
$ws.Range("N125").Value = -89499.336

$ws.Range("H134").Value = 1663184
$ws.Range("I134").Value = 1833100.9
$ws.Range("J134").Value = 6493.5
$ws.Range("K134").Value = 5499302.699999999
$ws.Range("L134").Value = 19480.5
$ws.Range("M134").Value = -5496767.699999999
$ws.Range("N134").Value = -24550.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 652.9091
$ws.Range("J22").Value = 833
$ws.Range("L22").Value = 833
$ws.Range("N22").Value = -1533

$ws.Range("H31").Value = 3361.1555
$ws.Range("I31").Value = 2061.4783
$ws.Range("J31").Value = 4719.909
$ws.Range("K31").Value = 2061.4783
$ws.Range("L31").Value = 4719.909
$ws.Range("M31").Value = -1766.4783
$ws.Range("N31").Value = -5309.909

$ws.Range("H34").Value = 3361.1555
$ws.Range("I34").Value = 2061.4783
$ws.Range("J34").Value = 4719.909
$ws.Range("K34").Value = 2061.4783
$ws.Range("L34").Value = 4719.909
$ws.Range("M34").Value = -1859.4783
$ws.Range("N34").Value = -5123.909

$ws.Range("H58").Value = 3175.7297
$ws.Range("I58").Value = 2989.4827
$ws.Range("K58").Value = 2989.4827
$ws.Range("M58").Value = -2786.4827

$ws.Range("H69").Value = 102887
$ws.Range("J69").Value = 107182.664
$ws.Range("L69").Value = 107182.664
$ws.Range("N69").Value = -108680.664

$ws.Range("H72").Value = 102887
$ws.Range("J72").Value = 107182.664
$ws.Range("L72").Value = 321547.992
$ws.Range("N72").Value = -329035.992

$ws.Range("H88").Value = 46496.5
$ws.Range("J88").Value = 46496.5
$ws.Range("L88").Value = 46496.5
$ws.Range("N88").Value = -47308.5

$ws.Range("H91").Value = 46496.5
$ws.Range("J91").Value = 46496.5
$ws.Range("L91").Value = 46496.5
$ws.Range("N91").Value = -49304.5

$ws.Range("H103").Value = 37677.812
$ws.Range("I103").Value = 9350.362999999999
$ws.Range("K103").Value = 9350.362999999999
$ws.Range("M103").Value = -8178.362999999999

$ws.Range("H105").Value = 2159.3572
$ws.Range("I105").Value = 1994.25
$ws.Range("K105").Value = 1994.25
$ws.Range("M105").Value = -247.25

$ws.Range("H107").Value = 1752.5
$ws.Range("I107").Value = 1063.5
$ws.Range("K107").Value = 1063.5
$ws.Range("M107").Value = 856.5

$ws.Range("H132").Value = 3253.554
$ws.Range("I132").Value = 3167.5789
$ws.Range("J132").Value = 3541.8235
$ws.Range("K132").Value = 9502.736699999999
$ws.Range("L132").Value = 10625.4705
$ws.Range("M132").Value = -6972.736699999999
$ws.Range("N132").Value = -15685.4705

$ws.Range("H134").Value = 2868.2856
$ws.Range("I134").Value = 2830.3076
$ws.Range("J134").Value = 2930
$ws.Range("K134").Value = 8490.9228
$ws.Range("L134").Value = 8790
$ws.Range("M134").Value = -5955.9228
$ws.Range("N134").Value = -13860

$ws.Range("H136").Value = 3175.7297
$ws.Range("I136").Value = 2989.4827
$ws.Range("K136").Value = 8968.4481
$ws.Range("M136").Value = -6418.4481

$ws.Range("H141").Value = 242690.92
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 11006.857
$ws.Range("J52").Value = 11006.857
$ws.Range("L52").Value = 33020.571
$ws.Range("N52").Value = -33552.571

$ws.Range("H107").Value = 925.96155
$ws.Range("I107").Value = 1306
$ws.Range("J107").Value = 757.05554
$ws.Range("K107").Value = 3918
$ws.Range("L107").Value = 2271.16662
$ws.Range("M107").Value = -1998
$ws.Range("N107").Value = -6111.16662

$ws.Range("H113").Value = 949.8570999999999
$ws.Range("J113").Value = 1173.1875
$ws.Range("L113").Value = 3519.5625
$ws.Range("N113").Value = -7859.5625

$ws.Range("H116").Value = 3934.75
$ws.Range("J116").Value = 5083.3335
$ws.Range("L116").Value = 15250.0005
$ws.Range("N116").Value = -22134.0005

$ws.Range("H117").Value = 222.83333
$ws.Range("I117").Value = 222.83333
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 668.49999
$ws.Range("L117").Value = 0
$ws.Range("M117").Value = 2773.50001
$ws.Range("N117").ClearContents()

$ws.Range("H130").Value = 3068.4285
$ws.Range("I130").Value = 2295.8
$ws.Range("K130").Value = 6887.400000000001
$ws.Range("M130").Value = -1867.400000000001

$ws.Range("H131").Value = 1536.3
$ws.Range("J131").Value = 1813.7407
$ws.Range("L131").Value = 5441.2221
$ws.Range("N131").Value = -15521.2221

$ws.Range("H136").Value = 1599.75
$ws.Range("I136").Value = 1599.75
$ws.Range("K136").Value = 4799.25
$ws.Range("M136").Value = 300.75

$ws.Range("H137").Value = 2707.7273
$ws.Range("I137").Value = 1465
$ws.Range("K137").Value = 4395
$ws.Range("M137").Value = 705

$ws.Range("H138").Value = 8194547.5
$ws.Range("I138").Value = 1921.2858
$ws.Range("J138").Value = 11212884
$ws.Range("K138").Value = 5763.857400000001
$ws.Range("L138").Value = 33638652
$ws.Range("M138").Value = -623.8574000000008
$ws.Range("N138").Value = -33648932

$ws.Range("H139").Value = 2570.2354
$ws.Range("I139").Value = 2266.2666
$ws.Range("K139").Value = 6798.7998
$ws.Range("M139").Value = -1658.7998

$ws.Range("H141").Value = 7500
$ws.Range("I141").Value = 5000
$ws.Range("K141").Value = 15000
$ws.Range("M141").Value = -9820

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 121.90476
$ws.Range("I2").Value = 63.8125
$ws.Range("J2").Value = 307.8
$ws.Range("K2").Value = 63.8125
$ws.Range("L2").Value = 307.8
$ws.Range("M2").Value = 49.1875
$ws.Range("N2").Value = -533.8

$ws.Range("H20").Value = 91665.664

$ws.Range("H70").Value = 4488.231
$ws.Range("I70").Value = 4399.6665
$ws.Range("K70").Value = 4399.6665
$ws.Range("M70").Value = -4129.6665

$ws.Range("H73").Value = 4488.231
$ws.Range("I73").Value = 4399.6665
$ws.Range("K73").Value = 4399.6665
$ws.Range("M73").Value = -3463.6665

$ws.Range("H97").Value = 849.2308
$ws.Range("I97").Value = 628.3333
$ws.Range("K97").Value = 628.3333
$ws.Range("M97").Value = -132.3333

$ws.Range("H100").Value = 36703.8
$ws.Range("J100").Value = 36703.8
$ws.Range("L100").Value = 36703.8
$ws.Range("N100").Value = -38867.8

$ws.Range("H102").Value = 1907.3334
$ws.Range("I102").Value = 1771.7142
$ws.Range("J102").Value = 2382
$ws.Range("K102").Value = 1771.7142
$ws.Range("L102").Value = 2382
$ws.Range("M102").Value = -149.7141999999999
$ws.Range("N102").Value = -5626

$ws.Range("H122").Value = 1581.6471
$ws.Range("I122").Value = 1198.1111
$ws.Range("K122").Value = 3594.3333
$ws.Range("M122").Value = -1144.3333

$ws.Range("H126").Value = 3591.3333
$ws.Range("I126").Value = 3449.6667
$ws.Range("K126").Value = 10349.0001
$ws.Range("M126").Value = -7879.000100000001

$ws.Range("H132").Value = 4275.4
$ws.Range("I132").Value = 4430.885
$ws.Range("K132").Value = 13292.655
$ws.Range("M132").Value = -10762.655

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6932.579
$ws.Range("J7").Value = 7299.9
$ws.Range("L7").Value = 7299.9
$ws.Range("N7").Value = -7523.9

$ws.Range("H16").Value = 851.82355
$ws.Range("I16").Value = 859.06665
$ws.Range("J16").Value = 797.5
$ws.Range("K16").Value = 859.06665
$ws.Range("L16").Value = 797.5
$ws.Range("M16").Value = -689.06665
$ws.Range("N16").Value = -1137.5

$ws.Range("H22").Value = 3485.9714
$ws.Range("I22").Value = 4635.8
$ws.Range("K22").Value = 4635.8
$ws.Range("M22").Value = -4340.8

$ws.Range("H27").Value = 3485.9714
$ws.Range("I27").Value = 4635.8
$ws.Range("K27").Value = 4635.8
$ws.Range("M27").Value = -4528.8

$ws.Range("H40").Value = 37040390
$ws.Range("I40").Value = 37040390
$ws.Range("K40").Value = 37040390
$ws.Range("M40").Value = -37040254

$ws.Range("H46").Value = 5582.7075
$ws.Range("I46").Value = 6849.227
$ws.Range("J46").Value = 4116.2104
$ws.Range("K46").Value = 6849.227
$ws.Range("L46").Value = 4116.2104
$ws.Range("M46").Value = -6661.227
$ws.Range("N46").Value = -4492.2104

$ws.Range("H55").Value = 556.63635
$ws.Range("I55").Value = 483.26666
$ws.Range("K55").Value = 483.26666
$ws.Range("M55").Value = -310.26666

$ws.Range("H100").Value = 2851.7273
$ws.Range("I100").Value = 2501.8235
$ws.Range("K100").Value = 2501.8235
$ws.Range("M100").Value = -1960.8235

$ws.Range("H120").Value = 111163
$ws.Range("J120").Value = 111163
$ws.Range("L120").Value = 111163
$ws.Range("N120").Value = -120839

$ws.Range("H122").Value = 16266.277
$ws.Range("I122").Value = 19258.2
$ws.Range("K122").Value = 57774.60000000001
$ws.Range("M122").Value = -55324.60000000001

$ws.Range("H126").Value = 6932.579
$ws.Range("J126").Value = 7299.9
$ws.Range("L126").Value = 21899.7
$ws.Range("N126").Value = -26839.7

$ws.Range("H132").Value = 716592.9399999999
$ws.Range("I132").Value = 835191.75
$ws.Range("K132").Value = 2505575.25
$ws.Range("M132").Value = -2503045.25

$ws.Range("H136").Value = 2705.5
$ws.Range("I136").Value = 2117.6
$ws.Range("J136").Value = 3881.3
$ws.Range("K136").Value = 6352.799999999999
$ws.Range("L136").Value = 11643.9
$ws.Range("M136").Value = -3802.799999999999
$ws.Range("N136").Value = -16743.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

$ws.Range("H62").Value = 3874.6667
$ws.Range("I62").Value = 2540.5
$ws.Range("J62").Value = 5208.8335
$ws.Range("K62").Value = 2540.5
$ws.Range("L62").Value = 5208.8335
$ws.Range("M62").Value = -1916.5
$ws.Range("N62").Value = -6456.8335

$ws.Range("H65").Value = 3874.6667
$ws.Range("I65").Value = 2540.5
$ws.Range("J65").Value = 5208.8335
$ws.Range("K65").Value = 12702.5
$ws.Range("L65").Value = 26044.1675
$ws.Range("M65").Value = -9582.5
$ws.Range("N65").Value = -32284.1675

$ws.Range("H70").Value = 74462.5
$ws.Range("I70").Value = 73000
$ws.Range("K70").Value = 73000
$ws.Range("M70").Value = -72685

$ws.Range("H73").Value = 74462.5
$ws.Range("I73").Value = 73000
$ws.Range("K73").Value = 73000
$ws.Range("M73").Value = -71908

$ws.Range("H75").Value = 78988
$ws.Range("J75").Value = 81247.5
$ws.Range("L75").Value = 81247.5
$ws.Range("N75").Value = -83119.5

$ws.Range("H78").Value = 78988
$ws.Range("J78").Value = 81247.5
$ws.Range("L78").Value = 243742.5
$ws.Range("N78").Value = -253102.5

$ws.Range("H100").Value = 1431.9286
$ws.Range("I100").Value = 1284.7
$ws.Range("K100").Value = 2569.4
$ws.Range("M100").Value = -2028.4

$ws.Range("H107").Value = 500.25
$ws.Range("I107").Value = 500.25
$ws.Range("K107").Value = 1500.75
$ws.Range("M107").Value = 419.25

$ws.Range("H122").Value = 31253604
$ws.Range("I122").Value = 35717100
$ws.Range("K122").Value = 107151300
$ws.Range("M122").Value = -107148850

$ws.Range("H130").Value = 79440.57000000001
$ws.Range("J130").Value = 79440.57000000001
$ws.Range("L130").Value = 79440.57000000001
$ws.Range("N130").Value = -89480.57000000001

$ws.Range("H132").Value = 39206.465
$ws.Range("I132").Value = 55804.844
$ws.Range("J132").Value = 4165.4443
$ws.Range("K132").Value = 167414.532
$ws.Range("L132").Value = 12496.3329
$ws.Range("M132").Value = -164884.532
$ws.Range("N132").Value = -17556.3329

$ws.Range("H136").Value = 19185.44
$ws.Range("I136").Value = 1145.1277
$ws.Range("J136").Value = 103974.9
$ws.Range("K136").Value = 3435.3831
$ws.Range("L136").Value = 311924.7
$ws.Range("M136").Value = -885.3831
$ws.Range("N136").Value = -317024.7
